# Update the "List of issues" tracker: record newly investigated
# countries (Rwanda, Equatorial Guinea, Burundi, Suriname, New
# Caledonia, Guyana) as additional rows under the existing table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append below the existing table (rows 1-23).
# Col A = Country, Col B = Problem described (re-using the existing
# "OSM issue with cables" note shared by several other countries).
$newRows = @(
    @{ Row = 24; Country = "Rwanda";           Height = 16 },
    @{ Row = 25; Country = "Equitorial Guinea"; Height = 32 },
    @{ Row = 26; Country = "Burundi ";          Height = 16 },
    @{ Row = 27; Country = "Suriname";          Height = 16 },
    @{ Row = 28; Country = "New Caledonia";     Height = 16 },
    @{ Row = 29; Country = "Guyana";            Height = 16 }
)

foreach ($entry in $newRows) {
    $r = $entry.Row

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = $entry.Country
    $cellA.WrapText = $true

    $cellB = $ws.Cells.Item($r, 2)
    $cellB.Value = "OSM issue with cables"
    $cellB.WrapText = $true

    $ws.Rows.Item($r).RowHeight = $entry.Height
}

# Update the selection / view to match where the author left off editing.
$ws.Range("A29").Select()
$excel.ActiveWindow.Zoom = 171
